# Refresh the cryptocurrency listing on Sheet1 (prices + 1h volume change)
# with the latest scrape results, and swap two pairs of rows that changed
# rank order (RenderToken<->MXToken, BabyDogeCoin<->Cronos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.823.75'
$ws.Range('E2').Value = '  -2.33%  '
# Row 3
$ws.Range('D3').Value = '1.561.31'
$ws.Range('E3').Value = '  -0.88%  '
# Row 4
$ws.Range('E4').Value = '  +0.08%  '
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '205.62'
$ws.Range('E5').Value = '  -0.89%  '
# Row 6
$ws.Range('E6').Value = '  -2.12%  '
# Row 7
$ws.Range('E7').Value = '  +0.11%  '
# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '21.66'
$ws.Range('E8').Value = '  -2.89%  '
# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.246'
$ws.Range('E9').Value = '  -1.27%  '
# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0583'
$ws.Range('E10').Value = '  -2.07%  '
# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0867'
$ws.Range('E11').Value = '  +0.08%  '
# Row 12
$ws.Range('D12').Value = '1.783.62'
$ws.Range('E12').Value = '  -0.77%  '
# Row 13
$ws.Range('D13').Value = '1.565.10'
$ws.Range('E13').Value = '  -0.65%  '
# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.72'
$ws.Range('E14').Value = '  -2.65%  '
# Row 15
$ws.Range('E15').Value = '  -1.70%  '
# Row 16
$ws.Range('D16').Value = '26.828.72'
$ws.Range('E16').Value = '  -2.28%  '
# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '60.97'
$ws.Range('E17').Value = '  -4.02%  '
# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.33'
$ws.Range('E18').Value = '  +0.33%  '
# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '213.25'
$ws.Range('E19').Value = '  -0.36%  '
# Row 20
$ws.Range('E20').Value = '  -2.30%  '
# Row 21
$ws.Range('E21').Value = '  +0.12%  '
# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  -1.07%  '
# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.18'
$ws.Range('E23').Value = '  -3.48%  '
# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.99'
$ws.Range('E24').Value = '  -0.57%  '
# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.57'
$ws.Range('E25').Value = '  +0.08%  '
# Row 26
$ws.Range('E26').Value = '  +0.12%  '
# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.90'
$ws.Range('E27').Value = '  -0.49%  '
# Row 28
$ws.Range('E28').Value = '  +0.08%  '
# Row 29
$ws.Range('E29').Value = '  -1.96%  '
# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0462'
$ws.Range('E30').Value = '  -2.17%  '
# Row 31
$ws.Range('E31').Value = '  -3.63%  '
# Row 32
$ws.Range('E32').Value = '  -1.41%  '
# Row 33
$ws.Range('D33').Value = '1.395.75'
$ws.Range('E33').Value = '  -0.15%  '
# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.90'
$ws.Range('E34').Value = '  -2.54%  '
# Row 35
$ws.Range('E35').Value = '  -3.10%  '
# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.27'
$ws.Range('E36').Value = '  -0.95%  '
# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.922'
$ws.Range('E37').Value = '  -1.84%  '
# Row 38
$ws.Range('E38').Value = '  -2.68%  '
# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.524'
$ws.Range('E39').Value = '  -1.27%  '
# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.809'
$ws.Range('E40').Value = '  -2.02%  '
# Row 41
$ws.Range('E41').Value = '  +0.08%  '
# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.992'
$ws.Range('E42').Value = '  -0.85%  '
# Row 43
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.18'
$ws.Range('E43').Value = '  -0.45%  '
# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.29'
$ws.Range('E44').Value = '  +0.44%  '
# Row 45
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.76'
$ws.Range('E45').Value = '  -3.14%  '
# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '62.78'
$ws.Range('E46').Value = '  -2.73%  '
# Row 47
$ws.Range('D47').Value = '1.696.76'
$ws.Range('E47').Value = '  -0.67%  '
# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '85.90'
$ws.Range('E48').Value = '  -0.14%  '
# Row 49
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0504'
$ws.Range('E49').Value = '  +1.89%  '
# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0981'
$ws.Range('E50').Value = '  -1.61%  '
# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0941'
$ws.Range('E51').Value = '  -1.36%  '
